$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L and M on this sheet are formatted as Text ("@"); writing a
# number straight into .Value on a Text-formatted cell gets coerced to a
# text string. Flip to General, write the number, then restore the
# original (Text) format so the stored style stays the same.
function Set-NumericValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# --- rows 237 & 245: daily new-case corrections (Oct 2020) ---
$ws.Range("C237").Value = 567
$ws.Range("C245").Value = 888

# --- rows 339-343: late Jan / early Feb 2021 update + newly reported day (row 343) ---
$ws.Range("C339").Value = 108
Set-NumericValue $ws.Range("L340") 2
$ws.Range("C341").Value = 54
$ws.Range("C342").Value = 91
Set-NumericValue $ws.Range("L342") 4

# Row 343 previously had no data reported for that day; fill in the new figures.
$ws.Range("C343").Value = 8
$ws.Range("E343").Value = 14
$ws.Range("F343").Value = 8
$ws.Range("G343").Value = 121
Set-NumericValue $ws.Range("L343") 0
Set-NumericValue $ws.Range("M343") 0
